$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$co1 = $ws.ChartObjects(1)
$co1.Left = 23.24984251968504
$co1.Top = 177.75
$co1.Width = 651.7500787401575
$co1.Height = 303.0

$co2 = $ws.ChartObjects(2)
$co2.Left = 620.9999212598425
$co2.Top = 166.5
$co2.Width = 529.5000787401575
$co2.Height = 348.7500787401575
